$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "SectionQuestionType" mini-table (column B) ---
$ws.Range("B9").Value = "SectionQuestionType"
$ws.Range("B11").Value = "QuestionAnswerSessionId"

# --- Extend the "QuestionAnswerSession" mini-table (column D) with two more rows ---
$ws.Range("D13").Value = "SectionId"
$ws.Range("D14").Value = "BookId"

# --- Remove the three mini-tables that used to live in columns F, H and J ---
$ws.Range("F9:F12").ClearContents()
$ws.Range("H9:H12").ClearContents()
$ws.Range("J9:J12").ClearContents()

# --- Column widths (best-fit recalculated by Excel after the content changes) ---
$ws.Columns.Item(2).ColumnWidth = 24.7109375
$ws.Columns.Item(10).ColumnWidth = 20.28515625

# --- Update the active selection to match the saved view state ---
$ws.Range("B12").Select()
